$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.041421333333334
$ws.Range("H2").Value = 6.124264
$ws.Range("I2").Value = 0.03013244547937793
$ws.Range("J2").Value = 0.03013244547937793
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.373057666666667
$ws.Range("N2").Value = 7.119173
$ws.Range("O2").Value = 0.03663411841135517
$ws.Range("P2").Value = 0.03663411841135517
$ws.Range("Q2").Value = 4.844410545963556
$ws.Range("R2").Value = 43.599694913672
$ws.Range("S2").Value = 0.001103875575715235
$ws.Range("T2").Value = 0.001103875575715235
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.041421333333334
$ws.Range("H3").Value = 6.124264
$ws.Range("I3").Value = 0.03013244547937793
$ws.Range("J3").Value = 0.03013244547937793
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 43.26393666666667
$ws.Range("N3").Value = 129.79181
$ws.Range("O3").Value = 0.6678877639810287
$ws.Range("P3").Value = 0.6678877639810287
$ws.Range("Q3").Value = 88.31992327531556
$ws.Range("R3").Value = 794.87930947784
$ws.Range("S3").Value = 0.02012509163450198
$ws.Range("T3").Value = 0.02012509163450198
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.041421333333334
$ws.Range("H4").Value = 6.124264
$ws.Range("I4").Value = 0.03013244547937793
$ws.Range("J4").Value = 0.03013244547937793
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4664283333333333
$ws.Range("N4").Value = 1.399285
$ws.Range("O4").Value = 0.007200495391983466
$ws.Range("P4").Value = 0.007200495391983466
$ws.Range("Q4").Value = 0.9521767501377778
$ws.Range("R4").Value = 8.56959075124
$ws.Range("S4").Value = 0.0002169685348234538
$ws.Range("T4").Value = 0.0002169685348234538
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.041421333333334
$ws.Range("H5").Value = 6.124264
$ws.Range("I5").Value = 0.03013244547937793
$ws.Range("J5").Value = 0.03013244547937793
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.67383333333333
$ws.Range("N5").Value = 56.0215
$ws.Range("O5").Value = 0.2882776222156329
$ws.Range("P5").Value = 0.2882776222156329
$ws.Range("Q5").Value = 38.12116174177778
$ws.Range("R5").Value = 343.090455676
$ws.Range("S5").Value = 0.008686509734337265
$ws.Range("T5").Value = 0.008686509734337265
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 34.65832566666667
$ws.Range("H6").Value = 103.974977
$ws.Range("I6").Value = 0.5115749950805638
$ws.Range("J6").Value = 0.5115749950805639
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.373057666666667
$ws.Range("N6").Value = 7.119173
$ws.Range("O6").Value = 0.03663411841135517
$ws.Range("P6").Value = 0.03663411841135517
$ws.Range("Q6").Value = 82.24620543711345
$ws.Range("R6").Value = 740.2158489340211
$ws.Range("S6").Value = 0.01874109894606981
$ws.Range("T6").Value = 0.01874109894606981
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 34.65832566666667
$ws.Range("H7").Value = 103.974977
$ws.Range("I7").Value = 0.5115749950805638
$ws.Range("J7").Value = 0.5115749950805639
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 43.26393666666667
$ws.Range("N7").Value = 129.79181
$ws.Range("O7").Value = 0.6678877639810287
$ws.Range("P7").Value = 0.6678877639810287
$ws.Range("Q7").Value = 1499.455606615375
$ws.Range("R7").Value = 13495.10045953837
$ws.Range("S7").Value = 0.3416746795729635
$ws.Range("T7").Value = 0.3416746795729636
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 34.65832566666667
$ws.Range("H8").Value = 103.974977
$ws.Range("I8").Value = 0.5115749950805638
$ws.Range("J8").Value = 0.5115749950805639
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.4664283333333333
$ws.Range("N8").Value = 1.399285
$ws.Range("O8").Value = 0.007200495391983466
$ws.Range("P8").Value = 0.007200495391983466
$ws.Range("Q8").Value = 16.16562507682722
$ws.Range("R8").Value = 145.490625691445
$ws.Range("S8").Value = 0.003683593394731564
$ws.Range("T8").Value = 0.003683593394731565
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 34.65832566666667
$ws.Range("H9").Value = 103.974977
$ws.Range("I9").Value = 0.5115749950805638
$ws.Range("J9").Value = 0.5115749950805639
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.67383333333333
$ws.Range("N9").Value = 56.0215
$ws.Range("O9").Value = 0.2882776222156329
$ws.Range("P9").Value = 0.2882776222156329
$ws.Range("Q9").Value = 647.2037971117223
$ws.Range("R9").Value = 5824.834174005501
$ws.Range("S9").Value = 0.147475623166799
$ws.Range("T9").Value = 0.147475623166799
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.322237
$ws.Range("H10").Value = 0.966711
$ws.Range("I10").Value = 0.004756386482002558
$ws.Range("J10").Value = 0.004756386482002558
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.373057666666667
$ws.Range("N10").Value = 7.119173
$ws.Range("O10").Value = 0.03663411841135517
$ws.Range("P10").Value = 0.03663411841135517
$ws.Range("Q10").Value = 0.7646869833336666
$ws.Range("R10").Value = 6.882182850003
$ws.Range("S10").Value = 0.0001742460255918507
$ws.Range("T10").Value = 0.0001742460255918507
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.322237
$ws.Range("H11").Value = 0.966711
$ws.Range("I11").Value = 0.004756386482002558
$ws.Range("J11").Value = 0.004756386482002558
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 43.26393666666667
$ws.Range("N11").Value = 129.79181
$ws.Range("O11").Value = 0.6678877639810287
$ws.Range("P11").Value = 0.6678877639810287
$ws.Range("Q11").Value = 13.94124115965667
$ws.Range("R11").Value = 125.47117043691
$ws.Range("S11").Value = 0.00317673233209428
$ws.Range("T11").Value = 0.00317673233209428
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.322237
$ws.Range("H12").Value = 0.966711
$ws.Range("I12").Value = 0.004756386482002558
$ws.Range("J12").Value = 0.004756386482002558
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.4664283333333333
$ws.Range("N12").Value = 1.399285
$ws.Range("O12").Value = 0.007200495391983466
$ws.Range("P12").Value = 0.007200495391983466
$ws.Range("Q12").Value = 0.1503004668483333
$ws.Range("R12").Value = 1.352704201635
$ws.Range("S12").Value = 0.00003424833894615187
$ws.Range("T12").Value = 0.00003424833894615187
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.322237
$ws.Range("H13").Value = 0.966711
$ws.Range("I13").Value = 0.004756386482002558
$ws.Range("J13").Value = 0.004756386482002558
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 18.67383333333333
$ws.Range("N13").Value = 56.0215
$ws.Range("O13").Value = 0.2882776222156329
$ws.Range("P13").Value = 0.2882776222156329
$ws.Range("Q13").Value = 6.017400031833334
$ws.Range("R13").Value = 54.1566002865
$ws.Range("S13").Value = 0.001371159785370276
$ws.Range("T13").Value = 0.001371159785370276
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 30.72629533333334
$ws.Range("H14").Value = 92.17888600000001
$ws.Range("I14").Value = 0.4535361729580556
$ws.Range("J14").Value = 0.4535361729580557
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.373057666666667
$ws.Range("N14").Value = 7.119173
$ws.Range("O14").Value = 0.03663411841135517
$ws.Range("P14").Value = 0.03663411841135517
$ws.Range("Q14").Value = 72.91527070903089
$ws.Range("R14").Value = 656.237436381278
$ws.Range("S14").Value = 0.01661489786397826
$ws.Range("T14").Value = 0.01661489786397827
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 30.72629533333334
$ws.Range("H15").Value = 92.17888600000001
$ws.Range("I15").Value = 0.4535361729580556
$ws.Range("J15").Value = 0.4535361729580557
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 43.26393666666667
$ws.Range("N15").Value = 129.79181
$ws.Range("O15").Value = 0.6678877639810287
$ws.Range("P15").Value = 0.6678877639810287
$ws.Range("Q15").Value = 1329.340495302629
$ws.Range("R15").Value = 11964.06445772366
$ws.Range("S15").Value = 0.3029112604414688
$ws.Range("T15").Value = 0.3029112604414689
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 30.72629533333334
$ws.Range("H16").Value = 92.17888600000001
$ws.Range("I16").Value = 0.4535361729580556
$ws.Range("J16").Value = 0.4535361729580557
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.4664283333333333
$ws.Range("N16").Value = 1.399285
$ws.Range("O16").Value = 0.007200495391983466
$ws.Range("P16").Value = 0.007200495391983466
$ws.Range("Q16").Value = 14.33161472183444
$ws.Range("R16").Value = 128.98453249651
$ws.Range("S16").Value = 0.003265685123482296
$ws.Range("T16").Value = 0.003265685123482296
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 30.72629533333334
$ws.Range("H17").Value = 92.17888600000001
$ws.Range("I17").Value = 0.4535361729580556
$ws.Range("J17").Value = 0.4535361729580557
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.67383333333333
$ws.Range("N17").Value = 56.0215
$ws.Range("O17").Value = 0.2882776222156329
$ws.Range("P17").Value = 0.2882776222156329
$ws.Range("Q17").Value = 573.7777180054445
$ws.Range("R17").Value = 5163.999462049001
$ws.Range("S17").Value = 0.1307443295291263
$ws.Range("T17").Value = 0.1307443295291263
